# Insert a new record (weekly "Poroto verde" price report) as row 52,
# shifting the existing rows 52:80 down to 53:81.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("52:52").Insert()

$ws.Range("A52").Value = 1
$ws.Range("B52").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C52").Value = "Arica y Parinacota"
$ws.Range("D52").Value = 44992
$ws.Range("E52").Value = 15
$ws.Range("F52").Value = 100112031
$ws.Range("G52").Value = "Poroto verde"
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 1300
$ws.Range("K52").Value = 900
$ws.Range("L52").Value = 1000
$ws.Range("M52").Value = 950
$ws.Range("N52").Value = "`$/kilo"
$ws.Range("O52").Value = "Región de Arica y Parinacota"
$ws.Range("P52").Value = 950
$ws.Range("Q52").Value = 1
$ws.Range("R52").Value = "Hortaliza"
